$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet is a daily/weekly price log for "Membrillo" (quince) at the
# Macroferia Regional de Talca. The update adds two new weekly price
# records, inserted in their date-sorted position within the existing
# block of rows (matching the row ordering convention already used in
# the sheet), pushing the following rows down.

# --- New record 1: insert at row 11 ---
$ws.Rows.Item(11).Insert()

$row11 = New-Object 'object[,]' 1,20
$row11[0,0]  = 5
$row11[0,1]  = "Macroferia Regional de Talca"
$row11[0,2]  = "Maule"
$row11[0,3]  = 44424
$row11[0,4]  = 7
$row11[0,5]  = "Fruta"
$row11[0,6]  = 100104
$row11[0,7]  = "Frutos de pepita"
$row11[0,8]  = 100104003
$row11[0,9]  = "Membrillo"
$row11[0,10] = "Champion"
$row11[0,11] = "Primera"
$row11[0,12] = 230
$row11[0,13] = 11000
$row11[0,14] = 11000
$row11[0,15] = 11000
$row11[0,16] = "`$/caja 18 kilos granel"
$row11[0,17] = "Región de O'Higgins"
$row11[0,18] = 611
$row11[0,19] = 18
$ws.Range("A11:T11").Value = $row11

# --- New record 2: insert at row 20 (after the shift caused above) ---
$ws.Rows.Item(20).Insert()

$row20 = New-Object 'object[,]' 1,20
$row20[0,0]  = 5
$row20[0,1]  = "Macroferia Regional de Talca"
$row20[0,2]  = "Maule"
$row20[0,3]  = 44427
$row20[0,4]  = 7
$row20[0,5]  = "Fruta"
$row20[0,6]  = 100104
$row20[0,7]  = "Frutos de pepita"
$row20[0,8]  = 100104003
$row20[0,9]  = "Membrillo"
$row20[0,10] = "Champion"
$row20[0,11] = "Primera"
$row20[0,12] = 200
$row20[0,13] = 11000
$row20[0,14] = 11000
$row20[0,15] = 11000
$row20[0,16] = "`$/caja 18 kilos granel"
$row20[0,17] = "Región de O'Higgins"
$row20[0,18] = 611
$row20[0,19] = 18
$ws.Range("A20:T20").Value = $row20

Write-Host "Inserted 2 rows; new dimension should be A1:T39"
